# Update "想去人数" (want-to-go count) values in column F across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1042  # was 1043
$ws.Range("F5").Value = 11  # was 10
$ws.Range("F7").Value = 530  # was 525
$ws.Range("F10").Value = 1888  # was 1885
$ws.Range("F11").Value = 4827  # was 4806
$ws.Range("F12").Value = 1369  # was 1365
$ws.Range("F14").Value = 2917  # was 2908
$ws.Range("F16").Value = 23  # was 22
$ws.Range("F17").Value = 1196  # was 1192
$ws.Range("F18").Value = 3964  # was 3953
$ws.Range("F19").Value = 920  # was 911
$ws.Range("F20").Value = 868  # was 863
$ws.Range("F21").Value = 1591  # was 1589
$ws.Range("F22").Value = 61  # was 59
$ws.Range("F23").Value = 2546  # was 2544
$ws.Range("F24").Value = 9  # was 8
$ws.Range("F31").Value = 1019  # was 1017
$ws.Range("F32").Value = 286  # was 284
$ws.Range("F33").Value = 71  # was 69
$ws.Range("F35").Value = 136  # was 132
$ws.Range("F36").Value = 1525  # was 1521
$ws.Range("F38").Value = 983  # was 982
$ws.Range("F40").Value = 217  # was 216
$ws.Range("F41").Value = 553  # was 550
$ws.Range("F42").Value = 183  # was 177
$ws.Range("F44").Value = 630  # was 629
$ws.Range("F45").Value = 347  # was 344
$ws.Range("F46").Value = 210  # was 207

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 629  # was 625

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 629  # was 625
$ws.Range("F3").Value = 1042  # was 1043
$ws.Range("F9").Value = 11  # was 10
$ws.Range("F11").Value = 530  # was 525
$ws.Range("F13").Value = 1888  # was 1885
$ws.Range("F14").Value = 4827  # was 4806
$ws.Range("F15").Value = 1369  # was 1365
$ws.Range("F18").Value = 2918  # was 2908
$ws.Range("F19").Value = 23  # was 22
$ws.Range("F20").Value = 1196  # was 1192
$ws.Range("F21").Value = 3964  # was 3953
$ws.Range("F22").Value = 920  # was 912
$ws.Range("F23").Value = 868  # was 864
$ws.Range("F24").Value = 1591  # was 1589
$ws.Range("F25").Value = 61  # was 59
$ws.Range("F26").Value = 2546  # was 2544
$ws.Range("F28").Value = 9  # was 8
$ws.Range("F36").Value = 1019  # was 1017
$ws.Range("F37").Value = 286  # was 284
$ws.Range("F38").Value = 1525  # was 1521
$ws.Range("F41").Value = 983  # was 982
$ws.Range("F44").Value = 553  # was 550
$ws.Range("F45").Value = 183  # was 177
$ws.Range("F46").Value = 630  # was 629
$ws.Range("F47").Value = 347  # was 344
$ws.Range("F48").Value = 210  # was 207
